$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in week 10 work: copy formatting/merge pattern from week 9 (rows 5-6)
# and set new content "Code chức năng" / "Hoàn thành" merged across rows 7-11.

$ws.Range("B5:C6").Copy() | Out-Null
$ws.Range("B7:C11").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("B7:B11").Merge() | Out-Null
$ws.Range("C7:C11").Merge() | Out-Null

$ws.Range("B7").Value = "Code chức năng"
$ws.Range("C7").Value = "Hoàn thành"

$ws.Range("B7:B11").Select() | Out-Null
